$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 62
$ws.Range("H62").Value = 285720860
$ws.Range("I62").Value = 333337660
$ws.Range("J62").Value = 250008260
$ws.Range("K62").Value = 333337660
$ws.Range("L62").Value = 250008260
$ws.Range("M62").Value = -333337036
$ws.Range("N62").Value = -250009508

# Row 65
$ws.Range("H65").Value = 285720860
$ws.Range("I65").Value = 333337660
$ws.Range("J65").Value = 250008260
$ws.Range("K65").Value = 1666688300
$ws.Range("L65").Value = 1250041300
$ws.Range("M65").Value = -1666685180
$ws.Range("N65").Value = -1250047540

# Row 96
$ws.Range("H96").Value = 1949888.4
$ws.Range("J96").Value = 3760232.8
$ws.Range("L96").Value = 11280698.4
$ws.Range("N96").Value = -11283444.4

# Row 99
$ws.Range("H99").Value = 1326.8
$ws.Range("I99").Value = 1326.8
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 3980.4
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -2482.4
$ws.Range("N99").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 5727.245
$ws.Range("I32").Value = 2170.0264
$ws.Range("J32").Value = 18015.818
$ws.Range("K32").Value = 2170.0264
$ws.Range("L32").Value = 18015.818
$ws.Range("M32").Value = -1883.0264
$ws.Range("N32").Value = -18589.818

# Row 61
$ws.Range("H61").Value = 1290.8541
$ws.Range("I61").Value = 1005.9091
$ws.Range("J61").Value = 1917.7333
$ws.Range("K61").Value = 1005.9091
$ws.Range("L61").Value = 1917.7333
$ws.Range("M61").Value = -793.9091
$ws.Range("N61").Value = -2341.7333

# Row 80
$ws.Range("H80").Value = 25275
$ws.Range("J80").Value = 25275
$ws.Range("L80").Value = 25275
$ws.Range("N80").Value = -27271

# Row 83
$ws.Range("H83").Value = 25275
$ws.Range("J83").Value = 25275
$ws.Range("L83").Value = 75825
$ws.Range("N83").Value = -85809

# Row 107
$ws.Range("H107").Value = 40000
$ws.Range("J107").Value = 40000
$ws.Range("L107").Value = 40000
$ws.Range("N107").Value = -47680

# Row 136
$ws.Range("H136").Value = 1290.8541
$ws.Range("I136").Value = 1005.9091
$ws.Range("J136").Value = 1917.7333
$ws.Range("K136").Value = 3017.7273
$ws.Range("L136").Value = 5753.199900000001
$ws.Range("M136").Value = -467.7273
$ws.Range("N136").Value = -10853.1999

$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 15387261
$ws.Range("I86").Value = 25002250
$ws.Range("J86").Value = 3280.2
$ws.Range("K86").Value = 25002250
$ws.Range("L86").Value = 3280.2
$ws.Range("M86").Value = -25001127
$ws.Range("N86").Value = -5526.2

# Row 89
$ws.Range("H89").Value = 15387261
$ws.Range("I89").Value = 25002250
$ws.Range("J89").Value = 3280.2
$ws.Range("K89").Value = 125011250
$ws.Range("L89").Value = 16401
$ws.Range("M89").Value = -125005634
$ws.Range("N89").Value = -27633

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 7274901.5
$ws.Range("J31").Value = 6063029.5
$ws.Range("L31").Value = 6063029.5
$ws.Range("N31").Value = -6063619.5

# Row 34
$ws.Range("H34").Value = 7274901.5
$ws.Range("J34").Value = 6063029.5
$ws.Range("L34").Value = 6063029.5
$ws.Range("N34").Value = -6063433.5

$ws = $wb.Worksheets.Item("CUL")
# Row 32
$ws.Range("H32").Value = 1530.7693
$ws.Range("I32").Value = 80
$ws.Range("J32").Value = 2437.5
$ws.Range("K32").Value = 240
$ws.Range("L32").Value = 7312.5
$ws.Range("M32").Value = 43
$ws.Range("N32").Value = -7878.5

# Row 39
$ws.Range("H39").Value = 2910.3157
$ws.Range("J39").Value = 2910.3157
$ws.Range("L39").Value = 8730.947100000001
$ws.Range("N39").Value = -9318.947100000001

# Row 68
$ws.Range("H68").Value = 1057.779
$ws.Range("I68").Value = 623.5263
$ws.Range("J68").Value = 1401.5625
$ws.Range("K68").Value = 1870.5789
$ws.Range("L68").Value = 4204.6875
$ws.Range("M68").Value = -1059.5789
$ws.Range("N68").Value = -5826.6875

# Row 71
$ws.Range("H71").Value = 1057.779
$ws.Range("I71").Value = 623.5263
$ws.Range("J71").Value = 1401.5625
$ws.Range("K71").Value = 5611.736699999999
$ws.Range("L71").Value = 12614.0625
$ws.Range("M71").Value = -1555.736699999999
$ws.Range("N71").Value = -20726.0625

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 22227110
$ws.Range("I70").Value = 36368564
$ws.Range("J70").Value = 4825.7144
$ws.Range("K70").Value = 36368564
$ws.Range("L70").Value = 4825.7144
$ws.Range("M70").Value = -36368294
$ws.Range("N70").Value = -5365.7144

# Row 73
$ws.Range("H73").Value = 22227110
$ws.Range("I73").Value = 36368564
$ws.Range("J73").Value = 4825.7144
$ws.Range("K73").Value = 36368564
$ws.Range("L73").Value = 4825.7144
$ws.Range("M73").Value = -36367628
$ws.Range("N73").Value = -6697.7144

# Row 108
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()

# Row 113
$ws.Range("H113").Value = 2087.795
$ws.Range("I113").Value = 1648.8667
$ws.Range("J113").Value = 2362.125
$ws.Range("K113").Value = 1648.8667
$ws.Range("L113").Value = 2362.125
$ws.Range("M113").Value = 521.1333
$ws.Range("N113").Value = -6702.125

# Row 114
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 1334165.9
$ws.Range("I22").Value = 3704146.5
$ws.Range("J22").Value = 1051.875
$ws.Range("K22").Value = 3704146.5
$ws.Range("L22").Value = 1051.875
$ws.Range("M22").Value = -3703851.5
$ws.Range("N22").Value = -1641.875

# Row 27
$ws.Range("H27").Value = 1334165.9
$ws.Range("I27").Value = 3704146.5
$ws.Range("J27").Value = 1051.875
$ws.Range("K27").Value = 3704146.5
$ws.Range("L27").Value = 1051.875
$ws.Range("M27").Value = -3704039.5
$ws.Range("N27").Value = -1265.875

# Row 110
$ws.Range("H110").Value = 35001
$ws.Range("J110").Value = 35001
$ws.Range("L110").Value = 35001
$ws.Range("N110").Value = -43181

$ws = $wb.Worksheets.Item("WVR")
# Row 70
$ws.Range("H70").Value = 18000
$ws.Range("J70").Value = 18000
$ws.Range("L70").Value = 18000
$ws.Range("N70").Value = -18630

# Row 73
$ws.Range("H73").Value = 18000
$ws.Range("J73").Value = 18000
$ws.Range("L73").Value = 18000
$ws.Range("N73").Value = -20184

# Row 132
$ws.Range("H132").Value = 2184.75
$ws.Range("I132").Value = 2048.5293
$ws.Range("J132").Value = 2515.5715
$ws.Range("K132").Value = 6145.5879
$ws.Range("L132").Value = 7546.7145
$ws.Range("M132").Value = -3615.5879
$ws.Range("N132").Value = -12606.7145
